# Apply updated odds values to the FlashScore weekly games sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (column -> new value)
$changes = @{
    2 = @{ 'K'=19; 'L'=1.13; 'M'=6; 'N'=1.44; 'O'=2.75; 'T'=13; 'U'=12 }
    4 = @{
        'G'=3.3; 'H'=2.7; 'I'=2.4; 'L'=1.36; 'P'=1.5; 'Q'=2.25; 'R'=1.7; 'S'=1.91;
        'T'=9; 'U'=18; 'V'=11; 'W'=50; 'X'=32; 'Y'=37; 'Z'=7.2; 'AA'=5.3; 'AB'=12.5;
        'AC'=60; 'AD'=500; 'AE'=7.1; 'AF'=11.75; 'AG'=9; 'AH'=28; 'AI'=21; 'AJ'=30
    }
    5 = @{ 'G'=5.5; 'H'=3.85; 'N'=1.62; 'V'=17; 'W'=110; 'X'=55; 'Y'=50; 'AA'=7.7; 'AJ'=21 }
    6 = @{
        'G'=1.42; 'H'=4.75; 'I'=7; 'L'=1.14; 'M'=5.5; 'N'=1.5; 'O'=2.5; 'W'=10;
        'Z'=17; 'AA'=9; 'AB'=17; 'AF'=41; 'AG'=21
    }
    7 = @{ 'G'=2.4; 'I'=2.75; 'T'=12; 'U'=15; 'V'=10; 'W'=23; 'AI'=19 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
